$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Resource Utilization" ---
$ws1 = $wb.Worksheets.Item("Resource Utilization")
$ws1.Range("B2").Value = 4.75
$ws1.Range("B3").Value = 1.77

# --- Sheet 2: "Activity Times" ---
$ws2 = $wb.Worksheets.Item("Activity Times")

# Capture current row 5 and row 6 values before overwriting, since rows 4-6
# get reshuffled (row4->row6 old content moved down, row5<->row6 swapped then row7 removed).
# Easiest: write the full desired final grid for rows 2-6 directly, then delete row 7.

$ws2.Range("A2").Value = "5.5.13 Real Property-Monthly Reviews-org"
$ws2.Range("B2").Value = "Process"
$ws2.Range("C2").Value = 10
$ws2.Range("D2").Value = 5
$ws2.Range("E2").Value = 21
$ws2.Range("F2").Value = 174
$ws2.Range("G2").Value = 67.2
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0

$ws2.Range("A3").Value = "Review AM using Asset Change Tracker (5.5.13.1)"
$ws2.Range("B3").Value = "Activity Step"
$ws2.Range("C3").Value = 8
$ws2.Range("D3").Value = 8
$ws2.Range("E3").Value = 7
$ws2.Range("F3").Value = 13
$ws2.Range("G3").Value = 9.75
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

$ws2.Range("A4").Value = "Note Accuracy in Asset Change Tracker (5.5.13.2)"
$ws2.Range("B4").Value = "Activity Step"
$ws2.Range("C4").Value = 6
$ws2.Range("D4").Value = 6
$ws2.Range("E4").Value = 3
$ws2.Range("F4").Value = 5
$ws2.Range("G4").Value = 4
$ws2.Range("H4").Value = 0
$ws2.Range("I4").Value = 0
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0

$ws2.Range("A5").Value = "Create/Post Journal Entries (5.5.13.4)"
$ws2.Range("B5").Value = "Stop"
$ws2.Range("C5").Value = 5
$ws2.Range("D5").Value = 5
$ws2.Range("E5").Value = 4
$ws2.Range("F5").Value = 6
$ws2.Range("G5").Value = 5.4
$ws2.Range("H5").Value = 0
$ws2.Range("I5").Value = 0
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0

$ws2.Range("A6").Value = "Work with REO RPO to Correct (5.5.13.3)"
$ws2.Range("B6").Value = "Activity Step"
$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 3
$ws2.Range("E6").Value = 80
$ws2.Range("F6").Value = 157
$ws2.Range("G6").Value = 110.67
$ws2.Range("H6").Value = 0
$ws2.Range("I6").Value = 0
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0

# Delete row 7 entirely (it no longer exists in the target)
$ws2.Rows.Item(7).Delete()

$wb.Save()
